# "20th April 1st update"
# Updates a few already-present 04/2020 case counts and inserts the first
# data points for 19/04/2020 and 20/04/2020 into the national time series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (14/04, 15/04, 16/04, 17/04) ---
$ws.Range("B30").Value = 1031
$ws.Range("B32").Value = 886
$ws.Range("B34").Value = 1061
$ws.Range("B36").Value = 922

# --- Insert 19/04/2020 right after 19/03/2020 (old row 39), before 20/03/2020 ---
$ws.Rows(40).Insert()
$ws.Range("A40").Value = "19/04/2020"
$ws.Range("B40").Value = 1580

# --- Insert 20/04/2020 right after 20/03/2020 (now row 41), before 21/03/2020 ---
$ws.Rows(42).Insert()
$ws.Range("A42").Value = "20/04/2020"
$ws.Range("B42").Value = 53
